$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new response row (row 5) by copying the formatting of the row
# above it (row 4) and then filling in the new values, mirroring how the
# previous rows of survey data were entered.
$ws.Range("A4:S4").Copy()
$ws.Range("A5:S5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A5").Value = 44203.741689814815
$ws.Range("B5").Value = "c"
$ws.Range("C5").Value = "3fmhl"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = 6
$ws.Range("H5").Value = 6
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 7
$ws.Range("L5").Value = 7
$ws.Range("M5").Value = 4
$ws.Range("N5").Value = 6
$ws.Range("O5").Value = 7
$ws.Range("P5").Value = 7
$ws.Range("Q5").Value = 6
$ws.Range("R5").Value = 6
$ws.Range("S5").Value = 4

# Match the row height used by the other data rows (16pt).
$ws.Rows.Item(5).RowHeight = 16

# Update the selected cell to reflect the author's cursor position after
# entering the new row of data.
$ws.Range("G10").Select() | Out-Null
